# Updates the "Price" (column D) and "Volume(1h)" (column E) figures on
# Sheet1 to the latest crypto snapshot, per the scheduled GitHub Actions
# refresh. Only the cells whose values actually changed are touched; all
# other cells (coin name, link, rank, unaffected price/volume rows) are
# left exactly as-is.
#
# Several "Price" values look numeric (e.g. 0.537, 35.80, 3.00) but must
# stay stored as literal text -- exactly like the source data -- so the
# trailing zeros / exact digit sequence survive instead of Excel silently
# reinterpreting the string as a number (which would drop trailing zeros,
# switch "." grouping, etc.). For those cells we briefly force a text
# number format while writing the value, then restore the cell's style
# to Normal so formatting ends up identical to the untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.791.25'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '2.930.26'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '376.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.575'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0851'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '3.400.10'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +66.74%  '
$ws.Range('D17').Value = '2.921.12'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.989'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '50.759.27'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.13%  '
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('D22').Value = '0.0₃0945'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.82%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.42'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('E31').Value = '  -3.65%  '
$ws.Range('E32').Value = '  +2.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.43'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '33.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0430'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('D48').Value = '1.998.81'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.258'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.67%  '
$ws.Range('E50').Value = '  -5.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.01%  '
